$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 309
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 3).Value = 46061
}
